# Scheduled-runner update: refresh Universalis market-price snapshots and
# recomputed profit columns (H:N) on the Leve profit tables for each class.
# Generated from the upstream price-refresh diff; CRP!M100 is removed (HQ
# item has no HQ recipe any more) and CUL!M107 is newly populated (item now
# has an HQ price), matching the upstream row edits.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 26683878
$ws.Range("I132").Value = 38613740
$ws.Range("J132").Value = 835844
$ws.Range("K132").Value = 115841220
$ws.Range("L132").Value = 2507532
$ws.Range("M132").Value = -115838690
$ws.Range("N132").Value = -2512592
$ws.Range("H137").Value = 613277.5
$ws.Range("I137").Value = 1539164.9
$ws.Range("J137").Value = 2585.8298
$ws.Range("K137").Value = 4617494.699999999
$ws.Range("L137").Value = 7757.4894
$ws.Range("M137").Value = -4614944.699999999
$ws.Range("N137").Value = -12857.4894

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4865.5425
$ws.Range("I32").Value = 4948.07
$ws.Range("K32").Value = 4948.07
$ws.Range("M32").Value = -4661.07
$ws.Range("H61").Value = 1426.6666
$ws.Range("I61").Value = 1419.5264
$ws.Range("J61").Value = 1494.5
$ws.Range("K61").Value = 1419.5264
$ws.Range("L61").Value = 1494.5
$ws.Range("M61").Value = -1207.5264
$ws.Range("N61").Value = -1918.5
$ws.Range("H94").Value = 34900
$ws.Range("J94").Value = 34900
$ws.Range("L94").Value = 34900
$ws.Range("N94").Value = -36702
$ws.Range("H132").Value = 2809.9524
$ws.Range("I132").Value = 1178.2727
$ws.Range("J132").Value = 4604.8
$ws.Range("K132").Value = 3534.8181
$ws.Range("L132").Value = 13814.4
$ws.Range("M132").Value = -1004.8181
$ws.Range("N132").Value = -18874.4
$ws.Range("H136").Value = 1426.6666
$ws.Range("I136").Value = 1419.5264
$ws.Range("J136").Value = 1494.5
$ws.Range("K136").Value = 4258.5792
$ws.Range("L136").Value = 4483.5
$ws.Range("M136").Value = -1708.5792
$ws.Range("N136").Value = -9583.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 26249
$ws.Range("J11").Value = 26249
$ws.Range("L11").Value = 26249
$ws.Range("N11").Value = -26529
$ws.Range("H50").Value = 31438.5
$ws.Range("J50").Value = 31438.5
$ws.Range("L50").Value = 31438.5
$ws.Range("N50").Value = -32688.5
$ws.Range("H51").Value = 23676.467
$ws.Range("J51").Value = 25389
$ws.Range("L51").Value = 25389
$ws.Range("N51").Value = -26861
$ws.Range("H60").Value = 25598.062
$ws.Range("J60").Value = 25598.062
$ws.Range("L60").Value = 25598.062
$ws.Range("N60").Value = -26620.062
$ws.Range("H61").Value = 23676.467
$ws.Range("J61").Value = 25389
$ws.Range("L61").Value = 25389
$ws.Range("N61").Value = -26085
$ws.Range("H74").Value = 31727.223
$ws.Range("J74").Value = 34532.5
$ws.Range("L74").Value = 34532.5
$ws.Range("N74").Value = -36280.5
$ws.Range("H77").Value = 31727.223
$ws.Range("J77").Value = 34532.5
$ws.Range("L77").Value = 103597.5
$ws.Range("N77").Value = -112333.5
$ws.Range("H99").Value = 14289410
$ws.Range("I99").Value = 33334814
$ws.Range("J99").Value = 5357.5
$ws.Range("K99").Value = 33334814
$ws.Range("L99").Value = 5357.5
$ws.Range("M99").Value = -33333316
$ws.Range("N99").Value = -8353.5
$ws.Range("H100").Value = 63333.332
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 63333.332
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 63333.332
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -65497.332
$ws.Range("H126").Value = 14289410
$ws.Range("I126").Value = 33334814
$ws.Range("J126").Value = 5357.5
$ws.Range("K126").Value = 100004442
$ws.Range("L126").Value = 16072.5
$ws.Range("M126").Value = -100001972
$ws.Range("N126").Value = -21012.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 284.33334
$ws.Range("I18").Value = 194.875
$ws.Range("K18").Value = 584.625
$ws.Range("M18").Value = -415.625
$ws.Range("H68").Value = 1303.4
$ws.Range("I68").Value = 971.1667
$ws.Range("J68").Value = 1445.7858
$ws.Range("K68").Value = 2913.5001
$ws.Range("L68").Value = 4337.357400000001
$ws.Range("M68").Value = -2102.5001
$ws.Range("N68").Value = -5959.357400000001
$ws.Range("H71").Value = 1303.4
$ws.Range("I71").Value = 971.1667
$ws.Range("J71").Value = 1445.7858
$ws.Range("K71").Value = 8740.5003
$ws.Range("L71").Value = 13012.0722
$ws.Range("M71").Value = -4684.5003
$ws.Range("N71").Value = -21124.0722
$ws.Range("H107").Value = 7158263.5
$ws.Range("I107").Value = 408.1316
$ws.Range("J107").Value = 15658217
$ws.Range("K107").Value = 1224.3948
$ws.Range("L107").Value = 46974651
$ws.Range("M107").Value = 695.6052
$ws.Range("N107").Value = -46978491
$ws.Range("H113").Value = 6945341
$ws.Range("I113").Value = 675.25
$ws.Range("K113").Value = 2025.75
$ws.Range("M113").Value = 144.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3457.48
$ws.Range("I126").Value = 2808.5
$ws.Range("J126").Value = 4836.5625
$ws.Range("K126").Value = 8425.5
$ws.Range("L126").Value = 14509.6875
$ws.Range("M126").Value = -5955.5
$ws.Range("N126").Value = -19449.6875
$ws.Range("H132").Value = 3501.25
$ws.Range("I132").Value = 1932.4546
$ws.Range("J132").Value = 4516.353
$ws.Range("K132").Value = 5797.3638
$ws.Range("L132").Value = 13549.059
$ws.Range("M132").Value = -3267.3638
$ws.Range("N132").Value = -18609.059
$ws.Range("H139").Value = 48999
$ws.Range("J139").Value = 48999
$ws.Range("L139").Value = 48999
$ws.Range("N139").Value = -59279

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3237.7046
$ws.Range("I132").Value = 2733.4167
$ws.Range("J132").Value = 3842.85
$ws.Range("K132").Value = 8200.250100000001
$ws.Range("L132").Value = 11528.55
$ws.Range("M132").Value = -5670.250100000001
$ws.Range("N132").Value = -16588.55
$ws.Range("H141").Value = 32275
$ws.Range("J141").Value = 32275
$ws.Range("L141").Value = 32275
$ws.Range("N141").Value = -42635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 28950
$ws.Range("J94").Value = 28950
$ws.Range("L94").Value = 28950
$ws.Range("N94").Value = -30752
$ws.Range("H126").Value = 563397.2
$ws.Range("I126").Value = 2063.2856
$ws.Range("J126").Value = 890842
$ws.Range("K126").Value = 6189.8568
$ws.Range("L126").Value = 2672526
$ws.Range("M126").Value = -3719.8568
$ws.Range("N126").Value = -2677466
